$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace exact text `oldText` with `newText` while keeping it an
# isolated run - i.e. prevent the engine's "merge adjacent runs that share
# identical rPr" behaviour from eating into the runs immediately before/after
# the match. We do this by briefly perturbing the formatting of the single
# character immediately before and after the match (so its rPr differs from
# the run being edited), performing the replace, then restoring the
# perturbed character's formatting back to its original value.
# ---------------------------------------------------------------------------
function Replace-Isolated {
    param($oldText, $newText)

    $rng = $d.Content
    $found = $rng.Find.Execute($oldText)
    $start = $rng.Start
    $end = $rng.End

    $hasBefore = $false
    $beforeSize = 0
    if ($start -gt 0) {
        $before = $d.Range($start - 1, $start)
        $bt = $before.Text
        if ($bt.Length -gt 0 -and [int][char]$bt[0] -ne 13 -and [int][char]$bt[0] -ne 7 -and [int][char]$bt[0] -ne 11) {
            $hasBefore = $true
            $beforeSize = $before.Font.Size
            $before.Font.Size = $beforeSize + 1
        }
    }

    $hasAfter = $false
    $afterSize = 0
    $docEnd = $d.Content.End
    if ($end -lt $docEnd) {
        $after = $d.Range($end, $end + 1)
        $at = $after.Text
        if ($at.Length -gt 0 -and [int][char]$at[0] -ne 13 -and [int][char]$at[0] -ne 7 -and [int][char]$at[0] -ne 11) {
            $hasAfter = $true
            $afterSize = $after.Font.Size
            $after.Font.Size = $afterSize + 1
        }
    }

    $rng2 = $d.Content
    $rng2.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

    $lenDiff = $newText.Length - $oldText.Length

    if ($hasBefore) {
        $before2 = $d.Range($start - 1, $start)
        $before2.Font.Size = $beforeSize
    }
    if ($hasAfter) {
        $after2 = $d.Range($end + $lenDiff, $end + $lenDiff + 1)
        $after2.Font.Size = $afterSize
    }
}

# ---------------------------------------------------------------------------
# Helper: replace exact text `oldText` (which may span multiple runs) with a
# single new run containing `newText` - this lets the engine's natural merge
# behaviour collapse the matched runs into one, which is what we want here.
# ---------------------------------------------------------------------------
function Replace-Merge {
    param($oldText, $newText)
    $rng = $d.Content
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# 1. Title
Replace-Isolated "The Enduring Enigma of Consciousness" "Biology: Exploring the Wonders of Life"

# 2. Author name: "Dr" + "." + " Anya Gupta" (3 runs) -> single run "Helen Robertson"
Replace-Merge "Dr. Anya Gupta" "Helen Robertson"

# 3. Email local-part/domain (keep ".org" runs untouched)
Replace-Isolated "agupta@neuro" "helenrobertson@edunet"

# 4. Body paragraph sentences (each remains its own isolated run)
Replace-Isolated "The intricate workings of consciousness have long captivated the curious minds of philosophers, scientists, and artists alike" "Our world teems with an astounding array of life, from the tiniest microbes to the majestic blue whales that grace our oceans"

Replace-Isolated " Its elusive nature has given rise to a plethora of theories and hypotheses, encompassing diverse fields such as neuroscience, psychology, and philosophy" " This captivating tapestry of living organisms, known as biodiversity, holds immense significance for our survival and well-being"

Replace-Isolated " Despite the advancements in modern science, consciousness remains an enigmatic phenomenon, evoking profound questions about our existence, perception, and relationship with the universe" " Biology, the study of life, unveils the intricate workings of these organisms, delving into their structure, function, growth, and evolution"

Replace-Isolated " Unraveling the complexities of consciousness presents an intellectual frontier that promises transformative insights into the essence of being" " This exploration guides us toward comprehending the beauty and complexity of the natural world"

Replace-Isolated "As we delve into the depths of consciousness, we encounter a myriad of perplexing questions" "Biology is a captivating subject that probes the mechanisms underlying the diversity of life"

Replace-Isolated " What are the physiological and neural correlates of consciousness? How do subjective experiences arise from electrochemical processes in the brain? What is the relationship between consciousness and the physical world we perceive? These questions challenge the boundaries of our understanding and invite us to explore the fundamental nature of reality" " This includes understanding the structure and function of cells, the basic unit of life, and how they work together to form tissues, organs, and organ systems in complex organisms"

Replace-Isolated " The journey to understanding consciousness promises to shed light on the interconnectedness of mind, body, and the cosmos, offering a glimpse into the deepest mysteries of human existence" " Biology also examines how organisms interact with their environment, revealing their fascinating adaptations and ecological relationships"

Replace-Isolated "Furthermore, the exploration of consciousness has profound implications for our understanding of free will, moral responsibility, and the nature of qualia" "From the grandeur of a towering forest ecosystem to the microscopic world of bacteria, biology unravels the interconnectedness of all living things"

Replace-Isolated " By unraveling the mechanisms underlying conscious experience, we may gain insights into the relationship between the conscious and unconscious mind, the genesis of creativity and inspiration, and the essence of self-awareness" " It elucidates the fundamental principles governing reproduction, genetics, and evolution, providing insights into the origin and diversification of life on Earth"

# 5. " The quest..." + "." + " It is a journey..." (3 runs) -> single run
Replace-Merge " The quest to comprehend consciousness is a testament to our insatiable curiosity, our desire to understand the universe and our place within it. It is a journey fraught with challenges, but the potential rewards are immense, promising transformative insights into the very essence of existence" " Understanding biology equips us with the knowledge to appreciate the fragility of our planet and the importance of preserving its biodiversity"

# 6. Summary paragraph
Replace-Isolated "The study of consciousness, an enduring enigma that has challenged scholars for centuries, presents a captivating frontier of intellectual exploration" "Biology, the study of life, unveils the intricate workings of living organisms, from the microscopic to the majestic"

# 7. " From the intricacies of " + lastRenderedPageBreak + "neural processes..." (2 runs) -> single run (drops the lastRenderedPageBreak)
Replace-Merge " From the intricacies of neural processes to the nature of subjective experience, the quest to understand consciousness promises transformative insights into our existence, perception, and relationship with the universe" " It delves into their structure, function, growth, and evolution, revealing the captivating tapestry of biodiversity that sustains our world"

# 8. " Its ramifications..." + "." + " As we delve deeper..." (3 runs) -> single run
Replace-Merge " Its ramifications extend to fundamental questions of free will, moral responsibility, and the essence of qualia. As we delve deeper into the mysteries of the conscious mind, we may unlock the secrets of creativity, inspiration, and self-awareness, gaining a profound understanding of the universe and our place within it" " Through its examination of cells, organisms, and their interactions with the environment, biology equips us with an understanding of the fundamental principles of life, guiding us toward appreciating the interconnectedness of all living things and the importance of preserving our planet's biodiversity"

# 9. Append an empty paragraph at the very end of the document body.
$d.Content.InsertParagraphAfter()
